$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.085.70"
$ws.Range("E2").Value = "  +3.20%  "
$ws.Range("D3").Value = "3.033.39"
$ws.Range("E3").Value = "  +1.89%  "
$ws.Range("E4").Value = "  +0.11%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "593.27"
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "154.00"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +8.15%  "
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("D8").Value = "3.030.99"
$ws.Range("E8").Value = "  +1.89%  "
$ws.Range("E9").Value = "  +0.30%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "6.88"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +14.69%  "
$ws.Range("E11").Value = "  +4.21%  "
$ws.Range("E12").Value = "  +1.94%  "
$ws.Range("E13").Value = "  +3.49%  "
$ws.Range("E14").Value = "  +5.09%  "
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").Value = "3.535.58"
$ws.Range("E16").Value = "  +1.97%  "
$ws.Range("D17").Value = "63.063.56"
$ws.Range("E17").Value = "  +3.23%  "
$ws.Range("E18").Value = "  +3.19%  "
$ws.Range("D19").Value = "3.033.54"
$ws.Range("E19").Value = "  +2.05%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "453.16"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.31%  "
$ws.Range("E21").Value = "  +1.35%  "
$ws.Range("E22").Value = "  +2.78%  "
$ws.Range("E23").Value = "  +3.31%  "
$ws.Range("E24").Value = "  +1.36%  "
$ws.Range("E25").Value = "  +9.84%  "
$ws.Range("E26").Value = "  +6.29%  "
$ws.Range("E27").Value = "  +4.47%  "
$ws.Range("E28").Value = "  -0.04%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "7.49"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +5.83%  "
$ws.Range("E30").Value = "  +11.06%  "
$ws.Range("E31").Value = "  +0.95%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.13%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "27.58"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.00%  "
$ws.Range("E34").Value = "  +2.36%  "
$ws.Range("D35").Value = "0.0₃0862"
$ws.Range("E35").Value = "  +6.37%  "
$ws.Range("E36").Value = "  +2.94%  "
$ws.Range("E37").Value = "  +3.17%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "3.15"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +11.00%  "
$ws.Range("E39").Value = "  +9.12%  "
$ws.Range("E40").Value = "  +2.81%  "
$ws.Range("E41").Value = "  +0.56%  "
$ws.Range("E42").Value = "  +1.53%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.309"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +16.48%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "43.93"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +12.32%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "392.34"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.64%  "
$ws.Range("E46").Value = "  +3.72%  "
$ws.Range("D47").Value = "2.721.16"
$ws.Range("E47").Value = "  +1.61%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "133.32"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +2.64%  "
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "25.34"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +9.30%  "
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "2.29"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +7.82%  "
